# Append 23 new departure rows (rows 249-271) for Friday, Jan 13
# to the "Main Data" sheet, mirroring the existing table layout
# (A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=TO, F=SHORT, G=AIRLINE,
#  H=MODEL, I=AIRCFAT ID, J=STATUS, K=(blank), L=DIFFERENCE, M=(blank)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 249 (NUMBER 248)
$ws.Cells.Item(249, 1).Value = 248
$ws.Cells.Item(249, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(249, 3).Value = "2:50 PM"
$ws.Cells.Item(249, 4).Value = "LO3838"
$ws.Cells.Item(249, 5).Value = "Warsaw"
$ws.Cells.Item(249, 6).Value = "(WAW)"
$ws.Cells.Item(249, 7).Value = "LOT "
$ws.Cells.Item(249, 8).Value = "E75S"
$ws.Cells.Item(249, 9).Value = "(SP-LIL)"
$ws.Cells.Item(249, 10).Value = "2:55 PM"
$ws.Cells.Item(249, 12).Value = "0 hours, 5 minutes"

# Row 250 (NUMBER 249)
$ws.Cells.Item(250, 1).Value = 249
$ws.Cells.Item(250, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(250, 3).Value = "3:15 PM"
$ws.Cells.Item(250, 4).Value = "FR3283"
$ws.Cells.Item(250, 5).Value = "Leeds"
$ws.Cells.Item(250, 6).Value = "(LBA)"
$ws.Cells.Item(250, 7).Value = "Ryanair "
$ws.Cells.Item(250, 8).Value = "B738"
$ws.Cells.Item(250, 9).Value = "(SP-RKM)"
$ws.Cells.Item(250, 10).Value = "3:23 PM"
$ws.Cells.Item(250, 12).Value = "0 hours, 8 minutes"

# Row 251 (NUMBER 250)
$ws.Cells.Item(251, 1).Value = 250
$ws.Cells.Item(251, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(251, 3).Value = "3:20 PM"
$ws.Cells.Item(251, 4).Value = "W61733"
$ws.Cells.Item(251, 5).Value = "Stockholm"
$ws.Cells.Item(251, 6).Value = "(NYO)"
$ws.Cells.Item(251, 7).Value = "Wizz Air "
$ws.Cells.Item(251, 8).Value = "A320"
$ws.Cells.Item(251, 9).Value = "(HA-LYM)"
$ws.Cells.Item(251, 10).Value = "3:26 PM"
$ws.Cells.Item(251, 12).Value = "0 hours, 6 minutes"

# Row 252 (NUMBER 251)
$ws.Cells.Item(252, 1).Value = 251
$ws.Cells.Item(252, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(252, 3).Value = "3:25 PM"
$ws.Cells.Item(252, 4).Value = "FR7150"
$ws.Cells.Item(252, 5).Value = "Krakow"
$ws.Cells.Item(252, 6).Value = "(KRK)"
$ws.Cells.Item(252, 7).Value = "Ryanair "
$ws.Cells.Item(252, 8).Value = "B738"
$ws.Cells.Item(252, 9).Value = "(SP-RSO)"
$ws.Cells.Item(252, 10).Value = "3:32 PM"
$ws.Cells.Item(252, 12).Value = "0 hours, 7 minutes"

# Row 253 (NUMBER 252)
$ws.Cells.Item(253, 1).Value = 252
$ws.Cells.Item(253, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(253, 3).Value = "3:45 PM"
$ws.Cells.Item(253, 4).Value = "SK760"
$ws.Cells.Item(253, 5).Value = "Copenhagen"
$ws.Cells.Item(253, 6).Value = "(CPH)"
$ws.Cells.Item(253, 7).Value = "SAS "
$ws.Cells.Item(253, 8).Value = "A20N"
$ws.Cells.Item(253, 9).Value = "(EI-SIE)"
$ws.Cells.Item(253, 10).Value = "3:49 PM"
$ws.Cells.Item(253, 12).Value = "0 hours, 4 minutes"

# Row 254 (NUMBER 253)
$ws.Cells.Item(254, 1).Value = 253
$ws.Cells.Item(254, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(254, 3).Value = "4:30 PM"
$ws.Cells.Item(254, 4).Value = "UNKNOWN"
$ws.Cells.Item(254, 5).Value = "Dresden"
$ws.Cells.Item(254, 6).Value = "(DRS)"
$ws.Cells.Item(254, 7).Value = "HeliService International "
$ws.Cells.Item(254, 8).Value = "A139"
$ws.Cells.Item(254, 9).Value = "(D-HHMH)"
$ws.Cells.Item(254, 10).Value = "4:51 PM"
$ws.Cells.Item(254, 12).Value = "0 hours, 21 minutes"

# Row 255 (NUMBER 254)
$ws.Cells.Item(255, 1).Value = 254
$ws.Cells.Item(255, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(255, 3).Value = "4:35 PM"
$ws.Cells.Item(255, 4).Value = "W61791"
$ws.Cells.Item(255, 5).Value = "Larnaca"
$ws.Cells.Item(255, 6).Value = "(LCA)"
$ws.Cells.Item(255, 7).Value = "Wizz Air "
$ws.Cells.Item(255, 8).Value = "A321"
$ws.Cells.Item(255, 9).Value = "(HA-LTB)"
$ws.Cells.Item(255, 10).Value = "4:43 PM"
$ws.Cells.Item(255, 12).Value = "0 hours, 8 minutes"

# Row 256 (NUMBER 255)
$ws.Cells.Item(256, 1).Value = 255
$ws.Cells.Item(256, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(256, 3).Value = "4:45 PM"
$ws.Cells.Item(256, 4).Value = "FR544"
$ws.Cells.Item(256, 5).Value = "London"
$ws.Cells.Item(256, 6).Value = "(STN)"
$ws.Cells.Item(256, 7).Value = "Ryanair "
$ws.Cells.Item(256, 8).Value = "B738"
$ws.Cells.Item(256, 9).Value = "(SP-RKQ)"
$ws.Cells.Item(256, 10).Value = "4:59 PM"
$ws.Cells.Item(256, 12).Value = "0 hours, 14 minutes"

# Row 257 (NUMBER 256)
$ws.Cells.Item(257, 1).Value = 256
$ws.Cells.Item(257, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(257, 3).Value = "5:15 PM"
$ws.Cells.Item(257, 4).Value = "LO3816"
$ws.Cells.Item(257, 5).Value = "Warsaw"
$ws.Cells.Item(257, 6).Value = "(WAW)"
$ws.Cells.Item(257, 7).Value = "LOT "
$ws.Cells.Item(257, 8).Value = "E195"
$ws.Cells.Item(257, 9).Value = "(SP-LNI)"
$ws.Cells.Item(257, 10).Value = "5:27 PM"
$ws.Cells.Item(257, 12).Value = "0 hours, 12 minutes"

# Row 258 (NUMBER 257)
$ws.Cells.Item(258, 1).Value = 257
$ws.Cells.Item(258, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(258, 3).Value = "5:35 PM"
$ws.Cells.Item(258, 4).Value = "BN1571"
$ws.Cells.Item(258, 5).Value = "La Roche-sur-Yon"
$ws.Cells.Item(258, 6).Value = "(EDM)"
$ws.Cells.Item(258, 7).Value = "Luxwing "
$ws.Cells.Item(258, 8).Value = "E55P"
$ws.Cells.Item(258, 9).Value = "(F-HLRS)"
$ws.Cells.Item(258, 10).Value = "5:53 PM"
$ws.Cells.Item(258, 12).Value = "0 hours, 18 minutes"

# Row 259 (NUMBER 258)
$ws.Cells.Item(259, 1).Value = 258
$ws.Cells.Item(259, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(259, 3).Value = "6:05 PM"
$ws.Cells.Item(259, 4).Value = "KL1922"
$ws.Cells.Item(259, 5).Value = "Amsterdam"
$ws.Cells.Item(259, 6).Value = "(AMS)"
$ws.Cells.Item(259, 7).Value = "KLM "
$ws.Cells.Item(259, 8).Value = "E190"
$ws.Cells.Item(259, 9).Value = "(PH-EZU)"
$ws.Cells.Item(259, 10).Value = "6:14 PM"
$ws.Cells.Item(259, 12).Value = "0 hours, 9 minutes"

# Row 260 (NUMBER 259)
$ws.Cells.Item(260, 1).Value = 259
$ws.Cells.Item(260, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(260, 3).Value = "6:10 PM"
$ws.Cells.Item(260, 4).Value = "W61705"
$ws.Cells.Item(260, 5).Value = "Barcelona"
$ws.Cells.Item(260, 6).Value = "(BCN)"
$ws.Cells.Item(260, 7).Value = "Wizz Air "
$ws.Cells.Item(260, 8).Value = "A21N"
$ws.Cells.Item(260, 9).Value = "(9H-WBU)"
$ws.Cells.Item(260, 10).Value = "6:16 PM"
$ws.Cells.Item(260, 12).Value = "0 hours, 6 minutes"

# Row 261 (NUMBER 260)
$ws.Cells.Item(261, 1).Value = 260
$ws.Cells.Item(261, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(261, 3).Value = "6:25 PM"
$ws.Cells.Item(261, 4).Value = "W61761"
$ws.Cells.Item(261, 5).Value = "Trondheim"
$ws.Cells.Item(261, 6).Value = "(TRD)"
$ws.Cells.Item(261, 7).Value = "Wizz Air "
$ws.Cells.Item(261, 8).Value = "A320"
$ws.Cells.Item(261, 9).Value = "(HA-LYO)"
$ws.Cells.Item(261, 10).Value = "6:29 PM"
$ws.Cells.Item(261, 12).Value = "0 hours, 4 minutes"

# Row 262 (NUMBER 261)
$ws.Cells.Item(262, 1).Value = 261
$ws.Cells.Item(262, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(262, 3).Value = "6:30 PM"
$ws.Cells.Item(262, 4).Value = "W61607"
$ws.Cells.Item(262, 5).Value = "London"
$ws.Cells.Item(262, 6).Value = "(LTN)"
$ws.Cells.Item(262, 7).Value = "Wizz Air "
$ws.Cells.Item(262, 8).Value = "A320"
$ws.Cells.Item(262, 9).Value = "(HA-LYS)"
$ws.Cells.Item(262, 10).Value = "6:32 PM"
$ws.Cells.Item(262, 12).Value = "0 hours, 2 minutes"

# Row 263 (NUMBER 262)
$ws.Cells.Item(263, 1).Value = 262
$ws.Cells.Item(263, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(263, 3).Value = "6:55 PM"
$ws.Cells.Item(263, 4).Value = "FR8510"
$ws.Cells.Item(263, 5).Value = "Oslo"
$ws.Cells.Item(263, 6).Value = "(TRF)"
$ws.Cells.Item(263, 7).Value = "Ryanair "
$ws.Cells.Item(263, 8).Value = "B738"
$ws.Cells.Item(263, 9).Value = "(SP-RSL)"
$ws.Cells.Item(263, 10).Value = "7:17 PM"
$ws.Cells.Item(263, 12).Value = "0 hours, 22 minutes"

# Row 264 (NUMBER 263)
$ws.Cells.Item(264, 1).Value = 263
$ws.Cells.Item(264, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(264, 3).Value = "7:00 PM"
$ws.Cells.Item(264, 4).Value = "W61681"
$ws.Cells.Item(264, 5).Value = "Milan"
$ws.Cells.Item(264, 6).Value = "(BGY)"
$ws.Cells.Item(264, 7).Value = "Wizz Air "
$ws.Cells.Item(264, 8).Value = "A320"
$ws.Cells.Item(264, 9).Value = "(HA-LYM)"
$ws.Cells.Item(264, 10).Value = "7:02 PM"
$ws.Cells.Item(264, 12).Value = "0 hours, 2 minutes"

# Row 265 (NUMBER 264)
$ws.Cells.Item(265, 1).Value = 264
$ws.Cells.Item(265, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(265, 3).Value = "7:05 PM"
$ws.Cells.Item(265, 4).Value = "W61615"
$ws.Cells.Item(265, 5).Value = "Leeds"
$ws.Cells.Item(265, 6).Value = "(LBA)"
$ws.Cells.Item(265, 7).Value = "Wizz Air "
$ws.Cells.Item(265, 8).Value = "A321"
$ws.Cells.Item(265, 9).Value = "(HA-LXL)"
$ws.Cells.Item(265, 10).Value = "7:10 PM"
$ws.Cells.Item(265, 12).Value = "0 hours, 5 minutes"

# Row 266 (NUMBER 265)
$ws.Cells.Item(266, 1).Value = 265
$ws.Cells.Item(266, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(266, 3).Value = "7:50 PM"
$ws.Cells.Item(266, 4).Value = "FR6111"
$ws.Cells.Item(266, 5).Value = "Lublin"
$ws.Cells.Item(266, 6).Value = "(LUZ)"
$ws.Cells.Item(266, 7).Value = "Ryanair "
$ws.Cells.Item(266, 8).Value = "B738"
$ws.Cells.Item(266, 9).Value = "(SP-RSW)"
$ws.Cells.Item(266, 10).Value = "8:01 PM"
$ws.Cells.Item(266, 12).Value = "0 hours, 11 minutes"

# Row 267 (NUMBER 266)
$ws.Cells.Item(267, 1).Value = 266
$ws.Cells.Item(267, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(267, 3).Value = "8:14 PM"
$ws.Cells.Item(267, 4).Value = "P81957"
$ws.Cells.Item(267, 5).Value = "Cologne"
$ws.Cells.Item(267, 6).Value = "(CGN)"
$ws.Cells.Item(267, 7).Value = "SprintAir "
$ws.Cells.Item(267, 8).Value = "SF34"
$ws.Cells.Item(267, 9).Value = "(SP-KPE)"
$ws.Cells.Item(267, 10).Value = "8:17 PM"
$ws.Cells.Item(267, 12).Value = "0 hours, 3 minutes"

# Row 268 (NUMBER 267)
$ws.Cells.Item(268, 1).Value = 267
$ws.Cells.Item(268, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(268, 3).Value = "8:30 PM"
$ws.Cells.Item(268, 4).Value = "3V4106"
$ws.Cells.Item(268, 5).Value = "Prague"
$ws.Cells.Item(268, 6).Value = "(PRG)"
$ws.Cells.Item(268, 7).Value = "ASL Airlines "
$ws.Cells.Item(268, 8).Value = "B738"
$ws.Cells.Item(268, 9).Value = "(OE-IXH)"
$ws.Cells.Item(268, 10).Value = "8:38 PM"
$ws.Cells.Item(268, 12).Value = "0 hours, 8 minutes"

# Row 269 (NUMBER 268)
$ws.Cells.Item(269, 1).Value = 268
$ws.Cells.Item(269, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(269, 3).Value = "8:40 PM"
$ws.Cells.Item(269, 4).Value = "FR6131"
$ws.Cells.Item(269, 5).Value = "Stockholm"
$ws.Cells.Item(269, 6).Value = "(ARN)"
$ws.Cells.Item(269, 7).Value = "Ryanair "
$ws.Cells.Item(269, 8).Value = "B738"
$ws.Cells.Item(269, 9).Value = "(SP-RKM)"
$ws.Cells.Item(269, 10).Value = "8:50 PM"
$ws.Cells.Item(269, 12).Value = "0 hours, 10 minutes"

# Row 270 (NUMBER 269)
$ws.Cells.Item(270, 1).Value = 269
$ws.Cells.Item(270, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(270, 3).Value = "9:20 PM"
$ws.Cells.Item(270, 4).Value = "LO3826"
$ws.Cells.Item(270, 5).Value = "Warsaw"
$ws.Cells.Item(270, 6).Value = "(WAW)"
$ws.Cells.Item(270, 7).Value = "LOT (Retro Livery) "
$ws.Cells.Item(270, 8).Value = "E75S"
$ws.Cells.Item(270, 9).Value = "(SP-LIM)"
$ws.Cells.Item(270, 10).Value = "9:55 PM"
$ws.Cells.Item(270, 12).Value = "0 hours, 35 minutes"

# Row 271 (NUMBER 270)
$ws.Cells.Item(271, 1).Value = 270
$ws.Cells.Item(271, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(271, 3).Value = "10:10 PM"
$ws.Cells.Item(271, 4).Value = "UNKNOWN"
$ws.Cells.Item(271, 5).Value = "Leipzig"
$ws.Cells.Item(271, 6).Value = "(LEJ)"
$ws.Cells.Item(271, 7).Value = "DHL "
$ws.Cells.Item(271, 8).Value = "B738"
$ws.Cells.Item(271, 9).Value = "(EC-IXO)"
$ws.Cells.Item(271, 10).Value = "10:14 PM"
$ws.Cells.Item(271, 12).Value = "0 hours, 4 minutes"
